{"js": "const results = context.document.body.search(\"-ok\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (const r of results.items) {\n  r.delete();\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"-ok\"\n$find.Replacement.Text = \"\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
